$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Determine the last used row in column C (falls back to 83 if detection fails)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 83 }

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (46075 -> 46076) for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46075) {
        $cell.Value2 = 46076
    }
}
